$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at E. This shifts old E -> F and old F -> G,
#    and Excel auto-adjusts all formulas/shared-formula refs that pointed at
#    those columns (matches the diff's column remap).
$ws.Columns("E").Insert()

# 2. New header row cells (E1 "Mean(Response^2)", H1 "Notes")
$ws.Range("E1").Value = "Mean(Response^2)"
$ws.Range("H1").Value = "Notes"

# 3. Row 8: back-transformation fixed to a manually computed static value
#    (was a formula "=(EXP(C8)-1)*100"); also add a note explaining why.
$ws.Range("F8").Value = -0.3191442
$ws.Range("H8").Value = "back transformation calculated in script to create Fig. 2"

# 4. Row 11: replace the "squared" back-transform inputs -- drop the old
#    MeanResponse (D11), add Mean(Response^2) value in E11, and rewrite the
#    percEstimate formula to use it.
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = 0.69601350890374436
$ws.Range("F11").Formula = "=((C11/2)/E11)*100"

# 5. Row 12: corrected MeanResponse input value
$ws.Range("D12").Value = 39.7098828

# 6. Row 13: keep the same numeric-display style on the (now-empty) E13,
#    matching the C13/D13 style (numFmtId 11, "0.00E+00").
$ws.Range("E13").NumberFormat = "0.00E+00"

# 7. Update the active-cell selection shown in the saved view.
$ws.Range("D13").Select()
